$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D407 (pupae count) on existing row 407
$ws.Range("D407").Value = 4

# Append new data rows 408-465 (vial, treatment, time_hours, pupae)
$data = @(
    ,@(1, "conditioned", 341, 0)
    ,@(1, "unconditioned", 341, 5)
    ,@(2, "conditioned", 341, 4)
    ,@(2, "unconditioned", 341, 0)
    ,@(3, "conditioned", 341, 0)
    ,@(3, "unconditioned", 341, 0)
    ,@(4, "conditioned", 341, 1)
    ,@(4, "unconditioned", 341, 1)
    ,@(5, "conditioned", 341, 0)
    ,@(5, "unconditioned", 341, 0)
    ,@(6, "conditioned", 341, 0)
    ,@(6, "unconditioned", 341, 1)
    ,@(7, "conditioned", 341, 1)
    ,@(7, "unconditioned", 341, 1)
    ,@(8, "conditioned", 341, 1)
    ,@(8, "unconditioned", 341, 0)
    ,@(9, "conditioned", 341, 3)
    ,@(9, "unconditioned", 341, 2)
    ,@(10, "conditioned", 341, 2)
    ,@(10, "unconditioned", 341, 0)
    ,@(11, "conditioned", 341, 1)
    ,@(11, "unconditioned", 341, 1)
    ,@(12, "conditioned", 341, 2)
    ,@(12, "unconditioned", 341, 1)
    ,@(13, "conditioned", 341, 0)
    ,@(13, "unconditioned", 341, 1)
    ,@(14, "conditioned", 341, 1)
    ,@(14, "unconditioned", 341, 0)
    ,@(15, "unconditioned", 341, 0)
    ,@(1, "conditioned", 356, 0)
    ,@(1, "unconditioned", 356, 3)
    ,@(2, "conditioned", 356, 0)
    ,@(2, "unconditioned", 356, 1)
    ,@(3, "conditioned", 356, 0)
    ,@(3, "unconditioned", 356, 2)
    ,@(4, "conditioned", 356, 0)
    ,@(4, "unconditioned", 356, 0)
    ,@(5, "conditioned", 356, 1)
    ,@(5, "unconditioned", 356, 0)
    ,@(6, "conditioned", 356, 0)
    ,@(6, "unconditioned", 356, 1)
    ,@(7, "conditioned", 356, 0)
    ,@(7, "unconditioned", 356, 0)
    ,@(8, "conditioned", 356, 1)
    ,@(8, "unconditioned", 356, 1)
    ,@(9, "conditioned", 356, 1)
    ,@(9, "unconditioned", 356, 1)
    ,@(10, "conditioned", 356, 0)
    ,@(10, "unconditioned", 356, 0)
    ,@(11, "conditioned", 356, 0)
    ,@(11, "unconditioned", 356, 1)
    ,@(12, "conditioned", 356, 1)
    ,@(12, "unconditioned", 356, 0)
    ,@(13, "conditioned", 356, 1)
    ,@(13, "unconditioned", 356, 1)
    ,@(14, "conditioned", 356, 1)
    ,@(14, "unconditioned", 356, 1)
    ,@(15, "unconditioned", 356, 1)
)

$startRow = 408
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$ws.Range("D466").Select()
